# Rename the "cruise_ID" attribute to "cruise" on the ColumnHeaders sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnHeaders")

$ws.Range("A5").Value = "cruise"

# Reflect the user's selection landing on the edited cell.
$ws.Range("A5").Select() | Out-Null
